# "Reverted to images on root" - strip the "data/" folder prefix from the
# image filenames in column A, and rename the header in B1 from "level" to
# "diff_level".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "diff_level"

$ws.Range("A2").Value = "img1.png"
$ws.Range("A3").Value = "img2.png"
$ws.Range("A4").Value = "img3.png"
$ws.Range("A5").Value = "img4.png"
$ws.Range("A6").Value = "img5.png"
$ws.Range("A7").Value = "img6.png"

# cosmetic: row heights grew slightly and the selection moved to D7 in the
# saved workbook
for ($r = 1; $r -le 7; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}

$ws.Range("D7").Select()
